# Fonctions : Add Candidate / Del Candidate / Add Elector (bug) / Dell Elector /
# Affichage intelligent des candidats
#
# Observed effect of the commit: a new candidate row was appended on the
# "candidat" sheet, and the active selection moved - ending with "candidat"
# as the active (selected) sheet, cell D7 selected there, while the
# "elector" sheet's selection moved to A5.

$wb = $excel.ActiveWorkbook

# --- Add Candidate -------------------------------------------------------
$wsCandidat = $wb.Worksheets.Item("candidat")

$wsCandidat.Range("A3").Value = "Lefou"
$wsCandidat.Range("B3").Value = "Oupas"
$wsCandidat.Range("C3").Value = "tg"
$wsCandidat.Range("D3").Value = "jesaispas"
$wsCandidat.Range("E3").Value = 3

# --- Update selections / active sheet ------------------------------------
# "elector" ends up no longer the active tab, selection parked on A5.
$wsElector = $wb.Worksheets.Item("elector")
$wsElector.Range("A5").Select()

# "candidat" becomes the active tab, selection on D7 (selected last).
$wsCandidat.Range("D7").Select()
